$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Margins (Project 1 Pars)")

# Revised margins/CI values for Project 1 Pars sheet
$ws.Range("D2").Value = 3.1647699999999999
$ws.Range("F2").Value = 3.1105700000000001
$ws.Range("G2").Value = 3.218969

$ws.Range("D3").Value = 3.2918829999999999
$ws.Range("F3").Value = 3.1932
$ws.Range("G3").Value = 3.3905660000000002

$ws.Range("D4").Value = 3.2128830000000002
$ws.Range("F4").Value = 3.0644960000000001
$ws.Range("G4").Value = 3.3612709999999999

$ws.Range("D5").Value = 3.3466939999999998
$ws.Range("F5").Value = 3.3110580000000001
$ws.Range("G5").Value = 3.3823300000000001

$ws.Range("D6").Value = 3.3808820000000002
$ws.Range("F6").Value = 3.3204020000000001
$ws.Range("G6").Value = 3.4413610000000001

$ws.Range("D7").Value = 3.369583
$ws.Range("F7").Value = 3.2792659999999998
$ws.Range("G7").Value = 3.4598990000000001

$ws.Range("D8").Value = 3.5286179999999998
$ws.Range("F8").Value = 3.5016750000000001
$ws.Range("G8").Value = 3.5555599999999998

$ws.Range("D9").Value = 3.4698799999999999
$ws.Range("F9").Value = 3.4311090000000002
$ws.Range("G9").Value = 3.508651

$ws.Range("D10").Value = 3.5262820000000001
$ws.Range("F10").Value = 3.4638239999999998
$ws.Range("G10").Value = 3.58874

$ws.Range("D11").Value = 3.7105419999999998
$ws.Range("F11").Value = 3.6744110000000001
$ws.Range("G11").Value = 3.7466729999999999

$ws.Range("D12").Value = 3.558878
$ws.Range("F12").Value = 3.501655
$ws.Range("G12").Value = 3.616101

$ws.Range("D13").Value = 3.6829809999999998
$ws.Range("F13").Value = 3.5868009999999999
$ws.Range("G13").Value = 3.7791619999999999

$ws.Range("D14").Value = 3.8924660000000002
$ws.Range("F14").Value = 3.837615
$ws.Range("G14").Value = 3.947317

$ws.Range("D15").Value = 3.6478769999999998
$ws.Range("F15").Value = 3.5531579999999998
$ws.Range("G15").Value = 3.7425959999999998

$ws.Range("D16").Value = 3.8396810000000001
$ws.Range("F16").Value = 3.6840980000000001
$ws.Range("G16").Value = 3.995263

# Match the author's final selection on this sheet (cell G19, just below the
# used range, as recorded in the saved sheetView).
$ws.Range("G19").Select()
